# Partial execution and adjustments
# Refresh the rpc-reply message-id UUIDs (and one commit flow-id) that were
# re-recorded for cells F2, H2 and I2 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: "First get config" response message-id
$ws.Range("F2").Value2 = $ws.Range("F2").Value2 -replace `
    "3315d179-330e-48da-8469-75da53992f20", "d3ac0de2-11e9-41b3-a901-c27618a0dfab"

# H2: "Edit config and commit" response - two message-ids plus the commit flow-id
$h2 = $ws.Range("H2").Value2
$h2 = $h2 -replace "88a79cff-c4ea-4b82-b7c4-0bb5fc362b53", "0f1f8bbe-807f-4941-bbae-c61609afa1df"
$h2 = $h2 -replace "a50d8f06-790b-4d23-9ef2-7fc13d7012a3", "dafe524b-4826-4880-bac2-ad0d6f0bd0a6"
$h2 = $h2 -replace 'nc-ext:flow-id="81"', 'nc-ext:flow-id="245"'
$ws.Range("H2").Value2 = $h2

# I2: "Second get config" response message-id
$ws.Range("I2").Value2 = $ws.Range("I2").Value2 -replace `
    "2bf056d0-b4bd-4b5d-9288-56095b0594da", "938e3aa1-eea7-46a9-be0a-07d64ea912a0"
